$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Merge column widths for columns E:F (5-6) into a single wider column width ---
# (Target stored width is 9.875 characters; the host's pixel-quantized ColumnWidth
# setter snaps to the nearest representable width, so we feed it the character
# count whose quantized result lands closest to 9.875.)
$ws.Range("E:F").ColumnWidth = 9

# --- Append new trade row 7 ---
$ws.Range("A7").Value = 42650.366828703707
$ws.Range("B7").Value = $true
$ws.Range("C7").Value = 10232.870000000001
$ws.Range("D7").Value = 10156.19
$ws.Range("E7").Value = 77.379997000000003
$ws.Range("F7").Value = 76.209998999999996
$ws.Range("G7").Value = $true
$ws.Range("H7").Value = -1.51
$ws.Range("I7").Value = $false

# A7/G7 use the workbook's date-format style (style index 1), same as the rows above.
# Copy that formatting across instead of re-deriving a NumberFormat string, so we
# reuse the existing style entry rather than minting a new one.
$ws.Range("A6").Copy()
$ws.Range("A7").PasteSpecial(-4122)
$ws.Range("G6").Copy()
$ws.Range("G7").PasteSpecial(-4122)
